$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve column D as text (prices are stored as text, e.g. "1.00", "0.110")
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "48.732.34"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.625.01"
$ws.Range("E3").Value = "  +3.89%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "110.16"
$ws.Range("E5").Value = "  +3.06%  "
$ws.Range("D6").Value = "322.19"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").Value = "0.520"
$ws.Range("E7").Value = "  -0.80%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").Value = "39.40"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").Value = "19.73"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "0.0808"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "7.19"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "3.020.57"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").Value = "2.623.86"
$ws.Range("E16").Value = "  +3.91%  "
$ws.Range("D17").Value = "0.857"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").Value = "48.649.18"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "12.81"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").Value = "6.66"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "0.0₃0939"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").Value = "270.57"
$ws.Range("E23").Value = "  -4.96%  "
$ws.Range("D24").Value = "68.98"
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").Value = "26.01"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").Value = "  +3.66%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "34.97"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "0.137"
$ws.Range("E31").Value = "  -4.45%  "
$ws.Range("D32").Value = "49.41"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "5.40"
$ws.Range("E33").Value = "  +1.79%  "
$ws.Range("D34").Value = "19.16"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "0.0792"
$ws.Range("E36").Value = "  +2.10%  "
$ws.Range("D37").Value = "4.91"
$ws.Range("E37").Value = "  +6.73%  "
$ws.Range("D38").Value = "2.02"
$ws.Range("E38").Value = "  +2.08%  "
$ws.Range("D39").Value = "3.12"
$ws.Range("E39").Value = "  +6.75%  "
$ws.Range("D40").Value = "124.35"
$ws.Range("E40").Value = "  +3.87%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "22.54"
$ws.Range("E41").Value = "  +3.28%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.110"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").Value = "2.13"
$ws.Range("E43").Value = "  -4.06%  "
$ws.Range("D44").Value = "0.0312"
$ws.Range("E44").Value = "  +2.69%  "
$ws.Range("D45").Value = "2.066.84"
$ws.Range("E45").Value = "  +3.15%  "
$ws.Range("D46").Value = "3.20"
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("E47").Value = "  +5.32%  "
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("D49").Value = "8.91"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").Value = "58.54"
$ws.Range("E50").Value = "  +3.66%  "
$ws.Range("D51").Value = "5.14"
$ws.Range("E51").Value = "  -0.82%  "
